$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 ("plotsFile"), shifting the existing
# rows 10-14 (dataFolder..outputFolder) down to rows 11-15.
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = "plotsFile"
$ws.Range("B10").Value = "Plots.xlsx"
$ws.Range("C10").Value = 'Name of the excel file with plot definitions. Must be located in the "paramsFolder"'

$ws.Range("B10").Select()
